$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160-190 down to 161-191
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new weekly price record
$ws.Cells.Item(160, 1).Value = 4
$ws.Cells.Item(160, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(160, 3).Value = "Los Lagos"
$ws.Cells.Item(160, 4).Value = 44476
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = 100114014
$ws.Cells.Item(160, 7).Value = "Betarraga"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 500
$ws.Cells.Item(160, 11).Value = 1000
$ws.Cells.Item(160, 12).Value = 1200
$ws.Cells.Item(160, 13).Value = 1100
$ws.Cells.Item(160, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 220
$ws.Cells.Item(160, 17).Value = 5
$ws.Cells.Item(160, 18).Value = "Hortaliza"
